# Auto-generated by diff-to-COM conversion
# Commit: Update automàtic: dades i banners [2026-02-23 04:20]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-23 04:18:15"
$ws.Range("E3").Value = "2026-02-23 04:18:17"
$ws.Range("E4").Value = "2026-02-23 04:18:20"
$ws.Range("H4").Value = "'92%"
$ws.Range("J4").Value = "1026.1 hPa"
$ws.Range("N4").Value = "3.7 °C 3:48 TU"
$ws.Range("O4").Value = "5.4 °C"
$ws.Range("E5").Value = "2026-02-23 04:18:22"
$ws.Range("N5").Value = "2.3 °C 3:41 TU"
$ws.Range("E6").Value = "2026-02-23 04:18:24"
$ws.Range("H6").Value = "'72%"
$ws.Range("J6").Value = "1025.6 hPa"
$ws.Range("N6").Value = "7.7 °C 3:59 TU"
$ws.Range("O6").Value = "9.6 °C"
$ws.Range("E7").Value = "2026-02-23 04:18:27"
$ws.Range("H7").Value = "'71%"
$ws.Range("J7").Value = "1025.2 hPa"
$ws.Range("N7").Value = "11.1 °C 3:48 TU"
$ws.Range("O7").Value = "12.1 °C"
$ws.Range("E8").Value = "2026-02-23 04:18:29"
$ws.Range("J8").Value = "1025.0 hPa"
$ws.Range("O8").Value = "13.2 °C"
$ws.Range("E9").Value = "2026-02-23 04:18:32"
$ws.Range("H9").Value = "'90%"
$ws.Range("N9").Value = "4.5 °C 3:42 TU"
$ws.Range("O9").Value = "7.2 °C"
$ws.Range("E10").Value = "2026-02-23 04:18:34"
$ws.Range("O10").Value = "4.5 °C"
$ws.Range("E11").Value = "2026-02-23 04:18:36"
$ws.Range("H11").Value = "'92%"
$ws.Range("N11").Value = "2.1 °C 3:51 TU"
$ws.Range("O11").Value = "3.1 °C"
$ws.Range("E12").Value = "2026-02-23 04:18:39"
$ws.Range("N12").Value = "4.1 °C 3:56 TU"
$ws.Range("O12").Value = "6.1 °C"
$ws.Range("E13").Value = "2026-02-23 04:18:41"
$ws.Range("H13").Value = "'90%"
$ws.Range("N13").Value = "-2.3 °C 3:54 TU"
$ws.Range("O13").Value = "-0.7 °C"
$ws.Range("E14").Value = "2026-02-23 04:18:43"
$ws.Range("H14").Value = "'84%"
$ws.Range("E15").Value = "2026-02-23 04:18:46"
$ws.Range("H15").Value = "'86%"
$ws.Range("N15").Value = "4.8 °C 3:55 TU"
$ws.Range("O15").Value = "7.2 °C"
$ws.Range("E16").Value = "2026-02-23 04:18:48"
$ws.Range("L16").Value = "28.4 km/h - 200º 3:51 TU"
$ws.Range("E17").Value = "2026-02-23 04:18:50"
$ws.Range("K17").Value = "-0.1 MJ/m2"
$ws.Range("O17").Value = "7.1 °C"
$ws.Range("E18").Value = "2026-02-23 04:18:53"
$ws.Range("J18").Value = "1026.2 hPa"
$ws.Range("N18").Value = "2.3 °C 3:51 TU"
$ws.Range("O18").Value = "3.3 °C"
$ws.Range("E19").Value = "2026-02-23 04:18:55"
$ws.Range("H19").Value = "'46%"
$ws.Range("E20").Value = "2026-02-23 04:18:58"
$ws.Range("H20").Value = "'39%"
$ws.Range("E21").Value = "2026-02-23 04:19:00"
$ws.Range("H21").Value = "'79%"
$ws.Range("N21").Value = "2.2 °C 3:33 TU"
$ws.Range("O21").Value = "4.0 °C"
$ws.Range("E22").Value = "2026-02-23 04:19:02"
$ws.Range("H22").Value = "'29%"
$ws.Range("N22").Value = "1.0 °C 3:59 TU"
$ws.Range("O22").Value = "2.2 °C"
$ws.Range("E23").Value = "2026-02-23 04:19:05"
$ws.Range("N23").Value = "1.1 °C 3:50 TU"
$ws.Range("O23").Value = "2.3 °C"
$ws.Range("E24").Value = "2026-02-23 04:19:07"
$ws.Range("J24").Value = "1027.6 hPa"
$ws.Range("N24").Value = "1.4 °C 3:52 TU"
$ws.Range("O24").Value = "3.0 °C"
$ws.Range("E25").Value = "2026-02-23 04:19:10"
$ws.Range("H25").Value = "'30%"
$ws.Range("N25").Value = "2.4 °C 3:49 TU"
$ws.Range("O25").Value = "3.3 °C"
$ws.Range("E26").Value = "2026-02-23 04:19:12"
$ws.Range("J26").Value = "1026.9 hPa"
$ws.Range("E27").Value = "2026-02-23 04:19:14"
$ws.Range("H27").Value = "'34%"
$ws.Range("O27").Value = "3.7 °C"
$ws.Range("E28").Value = "2026-02-23 04:19:17"
$ws.Range("J28").Value = "1027.5 hPa"
$ws.Range("N28").Value = "2.5 °C 3:56 TU"
$ws.Range("O28").Value = "4.0 °C"
$ws.Range("E29").Value = "2026-02-23 04:19:19"
$ws.Range("N29").Value = "3.5 °C 3:32 TU"
$ws.Range("E30").Value = "2026-02-23 04:19:21"
$ws.Range("J30").Value = "1025.8 hPa"
$ws.Range("N30").Value = "7.5 °C 3:59 TU"
$ws.Range("O30").Value = "8.3 °C"
$ws.Range("E31").Value = "2026-02-23 04:19:24"
$ws.Range("J31").Value = "1024.6 hPa"
$ws.Range("K31").Value = "-0.1 MJ/m2"
$ws.Range("N31").Value = "14.4 °C 3:32 TU"
$ws.Range("O31").Value = "15.0 °C"
$ws.Range("E32").Value = "2026-02-23 04:19:26"
$ws.Range("E33").Value = "2026-02-23 04:19:28"
$ws.Range("O33").Value = "2.8 °C"
$ws.Range("E34").Value = "2026-02-23 04:19:31"
$ws.Range("L34").Value = "14.4 km/h - 19º 3:41 TU"
$ws.Range("E35").Value = "2026-02-23 04:19:33"
$ws.Range("N35").Value = "9.5 °C 3:48 TU"
$ws.Range("O35").Value = "10.7 °C"
$ws.Range("E36").Value = "2026-02-23 04:19:36"
$ws.Range("J36").Value = "1025.5 hPa"
$ws.Range("N36").Value = "6.3 °C 3:35 TU"
$ws.Range("O36").Value = "7.1 °C"
$ws.Range("E37").Value = "2026-02-23 04:19:38"
$ws.Range("J37").Value = "1029.9 hPa"
$ws.Range("E38").Value = "2026-02-23 04:19:40"
$ws.Range("L38").Value = "14.4 km/h - 298º 3:59 TU"
$ws.Range("E39").Value = "2026-02-23 04:19:43"
$ws.Range("K39").Value = "-0.1 MJ/m2"
$ws.Range("O39").Value = "3.3 °C"
$ws.Range("E40").Value = "2026-02-23 04:19:45"
$ws.Range("H40").Value = "'91%"
$ws.Range("N40").Value = "1.3 °C 3:59 TU"
$ws.Range("O40").Value = "2.3 °C"
$ws.Range("E41").Value = "2026-02-23 04:19:48"
$ws.Range("H41").Value = "'84%"
$ws.Range("J41").Value = "1025.2 hPa"
$ws.Range("L41").Value = "7.6 km/h - 18º 3:37 TU"
$ws.Range("N41").Value = "6.4 °C 3:30 TU"
$ws.Range("O41").Value = "7.5 °C"
$ws.Range("E42").Value = "2026-02-23 04:19:50"
$ws.Range("N42").Value = "5.2 °C 3:36 TU"
$ws.Range("O42").Value = "6.2 °C"
$ws.Range("E43").Value = "2026-02-23 04:19:52"
$ws.Range("H43").Value = "'93%"
$ws.Range("N43").Value = "2.7 °C 3:59 TU"
$ws.Range("O43").Value = "4.4 °C"
$ws.Range("E44").Value = "2026-02-23 04:19:55"
$ws.Range("E45").Value = "2026-02-23 04:19:57"
$ws.Range("J45").Value = "1030.8 hPa"
$ws.Range("L45").Value = "13.7 km/h - 122º 3:48 TU"
$ws.Range("N45").Value = "2.2 °C 3:42 TU"
$ws.Range("O45").Value = "3.8 °C"
$ws.Range("E46").Value = "2026-02-23 04:19:59"
$ws.Range("J46").Value = "1027.5 hPa"
$ws.Range("N46").Value = "1.3 °C 3:48 TU"
$ws.Range("O46").Value = "2.5 °C"
